$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '29.483.51'
$cell.Style = $origStyle
$ws.Range('E2').Value = '  -1.42%  '

$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.850.69'
$cell.Style = $origStyle
$ws.Range('E3').Value = '  -0.57%  '

$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.9992'
$cell.Style = $origStyle

$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '243.06'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  -1.39%  '

$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.6586'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  +3.57%  '

$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.9999'
$cell.Style = $origStyle
$ws.Range('E7').Value = '  -0.05%  '

$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '48.00'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  +2.88%  '

$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.07504'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  +0.30%  '

$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.2991'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  -0.38%  '

$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '24.42'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  -0.78%  '

$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.07632'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  -0.68%  '

$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.841.31'
$cell.Style = $origStyle
$ws.Range('E13').Value = '  -1.03%  '

$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.019'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  -0.68%  '

$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.6858'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -0.57%  '

$ws.Range('E16').Value = '  -0.64%  '

$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.000009539'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  +1.48%  '

$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.140'
$cell.Style = $origStyle
$ws.Range('E18').Value = '  +0.94%  '

$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '29.519.78'
$cell.Style = $origStyle
$ws.Range('E19').Value = '  -1.16%  '

$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.085.74'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  -1.70%  '

$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '236.35'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  -1.24%  '

$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '12.59'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  -0.71%  '

$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.9997'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  -0.08%  '

$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.682'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  +4.33%  '

$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.1425'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  +0.35%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '156.94'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  -1.45%  '

$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.495'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -1.06%  '

$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '17.80'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -1.02%  '

$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.06028'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  -0.91%  '

$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.489'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  -1.12%  '

$ws.Range('E32').Value = '  -1.46%  '

$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.141'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  -0.08%  '

$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.074'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  -1.65%  '

$ws.Range('E35').Value = '  +1.39%  '

$ws.Range('E36').Value = '  -1.10%  '

$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.7230'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  -0.69%  '

$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.596'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  -0.82%  '

$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.802'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  -2.20%  '

$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.01781'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  -0.90%  '

$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.199.13'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  -1.94%  '

$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.240'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  -1.08%  '

$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.9073'
$cell.Style = $origStyle
$ws.Range('E43').Value = '  -2.03%  '

$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.9996'
$cell.Style = $origStyle
$ws.Range('E44').Value = '  -0.20%  '

$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.012.32'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -0.96%  '

$ws.Range('E46').Value = '  -0.47%  '

$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '66.28'
$cell.Style = $origStyle

$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.448'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +10.96%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.4058'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  -0.85%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.052'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -2.56%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.656'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +1.70%  '
